$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# 9x39 AP and DMG boost
# Row 21 = ammo_9x39_pab9 (Perf)
$ws.Range("G21").Value = 0.34
$ws.Range("H21").Value = 1.04

# Row 22 = ammo_9x39_ap (AP)
$ws.Range("G22").Value = 0.55
$ws.Range("H22").Value = 1.04

# Update the active cell selection as recorded in the sheet view
$ws.Range("H20").Select()
